$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Test_data" / 278 lookup pair from row 1 (E1:F1) down to row 2 (E2:F2)
$ws.Range("E1:F1").ClearContents()
$ws.Range("E2").Value = "Test_data"
$ws.Range("F2").Value = 278

# Update the % formulas in column C (rows 2-21) so they reference the
# relocated lookup cell $F$2 instead of $F$1
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("C$r").Formula = "=SUM(B$r/`$F`$2)"
}

# Move the active cell selection to L8 (matches the saved sheetView state)
$ws.Range("L8").Select()
